$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ternary")

# --- Column A: replace element list with new (sorted) dataset ---
$ws.Range("A2").Value = "Ba"
$ws.Range("A3").Value = "Ca"
$ws.Range("A4").Value = "Ce"
$ws.Range("A5").Value = "Co"
$ws.Range("A6").Value = "Eu"
$ws.Range("A7").Value = "Fe"
$ws.Range("A8").Value = "K"
$ws.Range("A9").Value = "La"
$ws.Range("A10").Value = "Na"
$ws.Range("A11").Value = "Nd"
$ws.Range("A12").Value = "Pr"
$ws.Range("A13").Value = "Sm"
$ws.Range("A14").Value = "Sr"
$ws.Range("A15").Value = "Yb"
$ws.Range("A16:A20").ClearContents()

# --- Column B updates ---
$ws.Range("B3").Value = "Os"
$ws.Range("B4").Value = "Ru"
$ws.Range("B5").Value = "Sn"
$ws.Range("B6").Value = "Tl"
$ws.Range("B7:B10").ClearContents()

# --- Column C updates ---
$ws.Range("C2").Value = "Sb"
$ws.Range("C3:C7").Clear()

# --- Remove now-unused trailing row ---
$ws.Rows.Item(22).Delete()

# --- View / selection state: Ternary becomes the active tab ---
$ws.Activate()
$ws.Range("F16").Select()

# --- Page setup (portrait) picked up on Ternary sheet ---
$ws.PageSetup.Orientation = 1
